$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-74 down to 55-75.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Cells.Item(54, 1).Value = 5
$ws.Cells.Item(54, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(54, 3).Value = "Maule"
$ws.Cells.Item(54, 4).Value = 44855
$ws.Cells.Item(54, 5).Value = 7
$ws.Cells.Item(54, 6).Value = 300000000
$ws.Cells.Item(54, 7).Value = "Espárragos"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 3000
$ws.Cells.Item(54, 11).Value = 1200
$ws.Cells.Item(54, 12).Value = 1200
$ws.Cells.Item(54, 13).Value = 1200
$ws.Cells.Item(54, 14).Value = "$/kilo"
$ws.Cells.Item(54, 15).Value = "Provincia de Linares"
$ws.Cells.Item(54, 16).Value = 1200
$ws.Cells.Item(54, 17).Value = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
